$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = "A deep dive into R Markdown"
$ws.Range("C10").Value = $true
$ws.Range("D11").Value = "Statistical learning: basics and linear regression"

$ws.Range("D12").Select()
